# Added larger data set to the pharma POS price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("(Aciclovir) 200 mg Tablet 1x") is no longer marked as sold.
$ws.Range("E2").Value = $false

# New row 5: another batch of "(Aciclovir) 200 mg Tablet 1x" - not sold.
$ws.Range("A5:E5").WrapText = $true
$ws.Range("D5").NumberFormat = "yyyy-mm-dd"
$ws.Range("A5").Value = $ws.Range("A2").Text
$ws.Range("B5").Value2 = 7
$ws.Range("C5").Value2 = 5
$ws.Range("D5").Value2 = $ws.Range("D2").Value2
$ws.Range("E5").Value = $false

# New row 6: another batch of "(Aciclovir) 200 mg Tablet 1x" - sold.
$ws.Range("A6:E6").WrapText = $true
$ws.Range("D6").NumberFormat = "yyyy-mm-dd"
$ws.Range("A6").Value = $ws.Range("A2").Text
$ws.Range("B6").Value2 = 2
$ws.Range("C6").Value2 = 0.6
$ws.Range("D6").Value2 = $ws.Range("D2").Value2
$ws.Range("E6").Value = $true
